$d = $word.ActiveDocument

# The bookmark "_GoBack" currently sits at the end of the "IP-Schutzklasse"
# paragraph. In the edited document it is relocated to the end of the new
# "LRCLK SYNC" paragraph, so remove it here and recreate it later.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Locate the end of the "...einhalten?" run (i.e. right after the text,
# before the paragraph mark) so the new list items can be inserted there.
$findRng = $d.Content
$findRng.Find.Execute("einhalten?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $findRng.End
$insertRng = $d.Range($insertPos, $insertPos)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wofür steht LRCLK SYNC</w:t></w:r><w:r><w:t>?</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wofür steht BCLK</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wofür steht DATA MISO:</w:t></w:r><w:r><w:t xml:space="preserve"> Daten MasterInputSlaveOutput </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertRng.InsertXML($xml)
